# Driver updated to catch SerialPort Exception and new Metadata tags added.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 9 (pushes existing rows 9-14 down to 10-15)
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with the new Raspberry Pi entry (AB-RPi02)
$ws.Range("A9").Value = "AB-RPi02"
$ws.Range("B9").Value = "Academic Building – Ground Floor"
$ws.Range("C9").Value = "192.168.136.58"
$ws.Range("D9").Value = "b8:27:eb:8b:bc:a4"
$ws.Range("E9").Value = "Deployed, not configured"

# The MB-RPI01 row (now row 15) status changes from "To be Deployed" to "Running"
$ws.Range("E15").Value = "Running"

# Update selection to match the author's final cursor position
[void]$ws.Range("E10").Select()

Write-Host "done"
